# Updated symbol list (Price / Volume(1h) columns) on Fri Feb 17 09:27:30 UTC 2023
# with GitHub Actions. Values are refreshed quotes for the crypto exchange
# tokens listed in the sheet; cells D/E must stay plain text (not numbers),
# matching the workbook's original inlineStr cell type, so each value is
# written with a leading apostrophe to force text, then the cell style is
# reset to "Normal" to strip Excel's auto-applied quote-prefix formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'309.92"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-3.63%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'49.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'2.16%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.121"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-2.59%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07754"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-4.28%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'4.535"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-1.17%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.377"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'14.78%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.570"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-4.37%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1219"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-6.37%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1981"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'1.72%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.04735"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'2.56%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.09312"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-1.47%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.1043"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.52%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001254"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-5.43%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.04173"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-2.72%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005803"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-2.09%"
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'2,020.74%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'-0.21%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'2.432"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'0.11%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'-0.24%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'8.012"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-1.31%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.1344"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-4.72%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.3033"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-2.98%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'-2.91%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.003934"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-7.43%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001350"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-0.04%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D38").Value = "'0.02599"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'-2.63%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.06273"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'11.54%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.01100"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'74.51%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007923"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'3.03%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1420"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-1.37%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.008371"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'8.74%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.008307"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'2.54%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3130"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-1.99%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00007341"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'4.95%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'-0.17%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.05321"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-0.42%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.002616"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-34.62%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002097"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.17%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0001997"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.17%"
$ws.Range("E51").Style = "Normal"